$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.059.99'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.876.51'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.00%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.52'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.91%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9991'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4915'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.82%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2925'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.60%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06614'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.14%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.882.11'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -0.70%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '16.55'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.73%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07213'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.6670'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.51%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '86.34'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.16%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.918'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.26%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.036.86'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000007810'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9991'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.125.33'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.91%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9969'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -0.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.791'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.863'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +3.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.146'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '152.36'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '142.89'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +8.75%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.16%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -3.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.385'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.215'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.14%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08792'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.78%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.998'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.57%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05074'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7219'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.01%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.659'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01875'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +13.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.687'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -3.17%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.81%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9302'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.33%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.792'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.24%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4241'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.55%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9983'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.17%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '103.25'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.384'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1280'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '32.87'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.11%  '
$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.3785'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.49%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.276'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.67%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.13%  '

Write-Host "Applied cryptos update"
